# "added 4wk low sales check"
#
# Refresh the forecast numbers on "Forecast Comparison" (MyForecast,
# Inventory Coverage, Stockout Risk, Reorder Urgency, Seasonality Index
# for weeks W10-W25) and roll the new totals up onto the "Summary" sheet
# (Total Forecast 16/8/4 Weeks, Max Forecast, Min Forecast).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. "Forecast Comparison" sheet - per-week forecast/inventory metrics
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Forecast Comparison")

# New values per row (only the columns that actually change are listed).
$forecastRows = @{
    2  = @{ D = 206; H = 9.42;               L = 0.8  }
    3  = @{ D = 207; H = 8.380000000000001;  L = 1.09 }
    4  = @{ D = 209; H = 7.31;               L = 1    }
    5  = @{ D = 210; H = 6.28;               L = 1.03 }
    6  = @{ D = 212; H = 5.23;               L = 0.88 }
    7  = @{ D = 213; H = 4.21;               L = 1.04 }
    8  = @{ D = 215; H = 3.18 }
    9  = @{ D = 216; H = 2.17;               L = 1.19 }
    10 = @{ D = 218; H = 1.16;               L = 1.04 }
    11 = @{ D = 219; H = 0.16; I = "High"; J = "Urgent"; L = 1.08 }
    12 = @{ D = 221; H = 0;    I = "High"; J = "Urgent"; L = 0.91 }
    13 = @{ D = 222; H = 0;    I = "High"; J = "Urgent"; L = 0.89 }
    14 = @{ D = 224; H = 0;    I = "High"; J = "Urgent"; L = 0.8100000000000001 }
    15 = @{ D = 225; H = 0;    I = "High"; J = "Urgent"; L = 1.15 }
    16 = @{ D = 227; H = 0;                              L = 1.17 }
    17 = @{ D = 228;                                      L = 0.98 }
}

foreach ($r in $forecastRows.Keys) {
    $row = $forecastRows[$r]
    if ($row.ContainsKey("D")) { $ws.Range("D$r").Value = $row.D }
    if ($row.ContainsKey("H")) { $ws.Range("H$r").Value = $row.H }
    if ($row.ContainsKey("I")) { $ws.Range("I$r").Value = $row.I }
    if ($row.ContainsKey("J")) { $ws.Range("J$r").Value = $row.J }
    if ($row.ContainsKey("L")) { $ws.Range("L$r").Value = $row.L }
}

# ---------------------------------------------------------------------
# 2. "Summary" sheet - roll-up totals derived from the new MyForecast
#    column above. These cells are stored as text, so force a text
#    number format before writing the new numeric-looking values,
#    otherwise Excel would auto-convert them to numbers.
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Summary")

$summaryRows = @{
    9  = "3472"   # Total Forecast (16 Weeks)
    10 = "1688"   # Total Forecast (8 Weeks)
    11 = "832"    # Total Forecast (4 Weeks)
    12 = "228"    # Max Forecast
    14 = "206"    # Min Forecast
}

foreach ($r in $summaryRows.Keys) {
    $cell = $ws2.Range("B$r")
    $cell.NumberFormat = "@"
    $cell.Value = $summaryRows[$r]
}
